$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d3Value = "Introduction
id=`"Par1`">In this paper, we describe our experience and early outcomes with critically unwell severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) patients who required extracorporeal membrane oxygenation (ECMO).

 We present our standard practices around ECMO decision-making, retrieval, cannulation, ventilation, anticoagulation, tracheostomy, imaging and steroids.


Methods
id=`"Par2`">A retrospective cohort study using data from the hospital notes on all SARS-CoV-2 patients who required extracorporeal support at St Bartholomew’s Hospital between 1 March 2020 and 31 July 2020. In total, this included 18 patients over this time period.


Results
id=`"Par3`">In total, 18 patients were managed with extracorporeal support and of these 14 survived (78%) with 4 deaths (22%).

 The mean duration from hospital admission to intubation was 4.1 ± 3.4 days, mean time from intubation to ECMO 2.3 ± 2 days and mean run on ECMO 17.7 ± 9.4 days.

 Survivor mean days from intubation to extubation was 20.6 ± 9.9 days and survivor mean days from intubation to tracheostomy decannulation 46.6 ± 15.3 days.

 Time from hospital admission to discharge in survivors was a mean of 57.2 ± 25.8 days.

 Of the patients requiring extracorporeal support, the initial mode was veno-venous (VV) in 15 (83%), veno-arterial (VA) in 2 (11%) and veno-venous-arterial (VVA) in 1 (6%).

 On VV extracorporeal support, 2 (11%) required additional VVA.

 Renal replacement therapy was required in 10 (56%) of the patients.

 Anticoagulation target anti-Xa of 0.2–0.4 was set, with 10 (56%) patients having a deep vein thrombosis or pulmonary embolism detected and 2 (11%) patients suffering an intracranial haemorrhage.

 Tracheostomy was performed in 9 (50%) of the patients and high-dose methylprednisolone was given to 7 (39%) of the patients.


Conclusion
id=`"Par4`">In our cohort of patients with severe SARS-CoV-2 respiratory failure, a long period of invasive ventilation and extracorporeal support was required but achieving good outcomes despite this.

 There was a significant burden of thromboembolic disease and renal injury.

 A significant proportion of patients required tracheostomy and steroids to facilitate weaning.


"
$e3Value = "[Waqas%Akhtar%NULL%1,   Olusegun%Olusanya%NULL%2,   Olusegun%Olusanya%NULL%0,   Marta Montero%Baladia%NULL%1,   Harriet%Young%NULL%1,   Sachin%Shah%s.shah4@nhs.net%1]"

$ws.Range("D3").Value2 = $d3Value
$ws.Range("E3").Value2 = $e3Value
